{"js": "// Fix systematic spacing issue between header bar and body text\n// - Collapse the three verbose CORE COMPETENCIES paragraphs into one summary line\n// - Add a new \"TECHNICAL SKILLS\" section (heading + three detail lines) near the end\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- Step 1: locate the three CORE COMPETENCIES detail paragraphs ----\nconst items = paragraphs.items;\n\nconst idxDataViz = items.findIndex(p =>\n  p.text.indexOf(\"Data Visualization & Design: Interactive Dashboards\") === 0\n);\nconst idxGeo = items.findIndex(p =>\n  p.text.indexOf(\"Geospatial Analysis & Mapping: Spatial Analysis\") === 0\n);\nconst idxTech = items.findIndex(p =>\n  p.text.indexOf(\"Technical Visualization: Programming\") === 0\n);\n\nif (idxDataViz === -1 || idxGeo === -1 || idxTech === -1) {\n  throw new Error(\"Could not locate CORE COMPETENCIES paragraphs\");\n}\n\n// Replace the first paragraph's text with the condensed summary line, then\n// delete the other two (now-redundant) paragraphs entirely.\nitems[idxDataViz].insertText(\n  \"Data Visualization & Design \u2022 Geospatial Analysis & Mapping \u2022 Technical Visualization\",\n  Word.InsertLocation.replace\n);\nitems[idxGeo].delete();\nitems[idxTech].delete();\n\nawait context.sync();\n\n// ---- Step 2: insert the new TECHNICAL SKILLS section ----\n// Re-fetch paragraphs since indices shifted after the deletes above.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst items2 = paragraphs2.items;\nconst idxLed = items2.findIndex(p =>\n  p.text.indexOf(\"Led multi-million dollar research projects\") !== -1\n);\nif (idxLed === -1) {\n  throw new Error(\"Could not locate the 'Led multi-million dollar...' paragraph\");\n}\n\nconst anchor = items2[idxLed];\n\n// Insert the three plain body lines (in reverse order) directly off the\n// Normal-styled anchor paragraph so they naturally inherit the Normal style\n// (no explicit style assignment needed -> matches the target markup, which\n// carries no <w:pPr> at all on these paragraphs).\nconst line3 = anchor.insertParagraph(\n  \"TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing\",\n  Word.InsertLocation.after\n);\nconst line2 = anchor.insertParagraph(\n  \"GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing\",\n  Word.InsertLocation.after\n);\nconst line1 = anchor.insertParagraph(\n  \"DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design\",\n  Word.InsertLocation.after\n);\n\n// Finally add the \"TECHNICAL SKILLS\" heading right after the anchor (pushing\n// the three lines above further down) and mark it as Heading 2.\nconst heading = anchor.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.after);\nheading.styleBuiltIn = Word.Style.heading2;\n\nawait context.sync();\n", "ps1": "# Fix systematic spacing issue between header bar and body text\n# - Collapse the three verbose CORE COMPETENCIES paragraphs into one summary line\n# - Add a new \"TECHNICAL SKILLS\" section (heading + three detail lines) near the end\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n\nfunction Find-ParagraphByPrefix($doc, $prefix) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# ---- Step 1: locate + collapse the three CORE COMPETENCIES detail paragraphs ----\n# NOTE: paragraph references in this COM model are position-based, so once a\n# paragraph is deleted, any *other* previously-fetched paragraph reference can\n# end up pointing at whatever now occupies its old slot. To stay safe, we\n# re-locate each paragraph by its text right before we touch it, rather than\n# collecting all three references up front.\n\n$dataVizPrefix = \"Data Visualization & Design: Interactive Dashboards\"\n$pDataViz = Find-ParagraphByPrefix $d $dataVizPrefix\nif ($pDataViz -eq $null) {\n    throw \"Could not locate the 'Data Visualization & Design' CORE COMPETENCIES paragraph\"\n}\n$condensed = \"Data Visualization & Design \" + $bullet + \" Geospatial Analysis & Mapping \" + $bullet + \" Technical Visualization\"\n$pDataViz.Range.Text = $condensed\n\n$geoPrefix = \"Geospatial Analysis & Mapping: Spatial Analysis\"\n$pGeo = Find-ParagraphByPrefix $d $geoPrefix\nif ($pGeo -eq $null) {\n    throw \"Could not locate the 'Geospatial Analysis & Mapping' CORE COMPETENCIES paragraph\"\n}\n$pGeo.Range.Delete()\n\n$techPrefix = \"Technical Visualization: Programming\"\n$pTech = Find-ParagraphByPrefix $d $techPrefix\nif ($pTech -eq $null) {\n    throw \"Could not locate the 'Technical Visualization' CORE COMPETENCIES paragraph\"\n}\n$pTech.Range.Delete()\n\n# ---- Step 2: insert the new TECHNICAL SKILLS section ----\n$ledPrefix = $bullet + \" Led multi-million dollar research projects\"\n$pLed = Find-ParagraphByPrefix $d $ledPrefix\nif ($pLed -eq $null) {\n    throw \"Could not locate the 'Led multi-million dollar...' paragraph\"\n}\n\n# Each InsertParagraphAfter() call on $pLed.Range inserts immediately after\n# $pLed (pushing any previously-inserted paragraph further down), so build\n# the new paragraphs in reverse order to end up with the right final order:\n# TECHNICAL SKILLS (heading), then the three detail lines.\n$pLed.Range.InsertParagraphAfter()\n$line3Para = $pLed.Next()\n$line3Para.Range.Text = \"TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing\"\n\n$pLed.Range.InsertParagraphAfter()\n$line2Para = $pLed.Next()\n$line2Para.Range.Text = \"GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing\"\n\n$pLed.Range.InsertParagraphAfter()\n$line1Para = $pLed.Next()\n$line1Para.Range.Text = \"DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design\"\n\n$pLed.Range.InsertParagraphAfter()\n$headingPara = $pLed.Next()\n$headingPara.Range.Text = \"TECHNICAL SKILLS\"\n$headingPara.Style = \"Heading 2\"\n"}
